$wb = $excel.ActiveWorkbook

# 1) Fill the "Procesos" sheet with its new data (A1:A20) and move the
#    selection to A21 - this mirrors what the author typed before adding
#    the new sheet, so it must happen first (selecting a range would
#    otherwise make "Procesos" the active/last-selected sheet).
$procesos = $wb.Worksheets.Item("Procesos")
$valores = @("a","d","a","a","d","c","b","d","b","c","b","a","a","c","a","b","a","d","b","d")
for ($i = 0; $i -lt $valores.Length; $i++) {
    $procesos.Cells.Item($i + 1, 1).Value = $valores[$i]
}
$procesos.Range("A21").Select() | Out-Null

# 2) Append the new "Integración" worksheet after the last existing sheet
#    (Excel's Worksheets.Add defaults to inserting before the active
#    sheet, so we explicitly target the position after the current last
#    sheet).
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$nuevaHoja = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$nuevaHoja.Name = "Integración"

# Adding a sheet automatically activates it, which is what we want: it
# becomes the selected tab (activeTab moves to its index, tabSelected
# moves off "Marco Conceptual" and off "Procesos").
$nuevaHoja.Activate() | Out-Null
